# Regenerate the handoff report: swap the old file GUID for a new one and
# bump the "Ready for handoff" / "Latest Handoff" timestamps.
#
# Old source file id : 1730f04a-5189-45e5-9e08-78893916f85e
# New source file id : 3a009bb0-b78e-4caf-9503-e1452583d4ef
#
# Old target hash (zh-cn/de-de xliff) : f53991ab623be33e4a4fdba43b9bee789468f759
# New target hash (zh-cn/de-de xliff) : 1d74d3ffb5882d5aeba39c76dfb0a896795094c7

$wb = $excel.ActiveWorkbook

$oldId = "1730f04a-5189-45e5-9e08-78893916f85e"
$newId = "3a009bb0-b78e-4caf-9503-e1452583d4ef"

$oldHash = "f53991ab623be33e4a4fdba43b9bee789468f759"
$newHash = "1d74d3ffb5882d5aeba39c76dfb0a896795094c7"

# The hyperlink Address (it points at the blob on the old commit SHA and
# keeps referencing the old file id - that part of the URL is untouched
# by this commit, only the cell's displayed text changes) is identical
# on all three sheets.
$addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/433206ab4cbe0672011d5390f7f226e61d64f056/e2e/$oldId.md"

# ---------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newId.md"
$ws.Range("B2").Value = "e2e\$newId.md"
$ws.Range("G2").Value = "2016-08-29 19:03:32"

# This engine's COM layer can't read back an existing Hyperlink's
# properties (only write them), and re-writing TextToDisplay/Address on
# an existing Hyperlinks item appends a duplicate instead of editing it
# in place - so replace it: clear the sheet's hyperlinks, then recreate
# the single link with the same Address and the refreshed display text.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $addr, "", "", "e2e\$newId.md")

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newId.md"
$ws.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-29 19:03:28"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addr, "", "", "$newId.md")

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newId.md"
$ws.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-29 19:03:32"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $addr, "", "", "$newId.md")
